# "adding simple report and tweaking data"
#
# The underlying date/value data on the "Data1" sheet is unchanged; only the
# two coverage-column headers in row 3 are renamed, and the active selection
# on that sheet moves from E26 to C3 (as if the author scrolled back up to
# review/tweak the header row after adding it).

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data1")

# Rename the header labels in row 3 (B3/C3) to the new coverage columns.
$wsData.Range("B3").Value = "GF coverage"
$wsData.Range("C3").Value = "BF coverage"

# Reflect the author's final selection being back up on the header row.
$wsData.Activate()
$wsData.Range("C3").Select()
